$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume/name/link values cell by cell.
# Cells whose new text would otherwise be auto-recognized by Excel as a
# number (e.g. "1.00", "0.479") are forced to Text format first, then the
# cell style is reset back to "Normal" so no stray number-format style is
# left attached to the cell.

$ws.Range("D2").Value = '61.112.32'
$ws.Range("E2").Value = '  -1.75%  '
$ws.Range("D3").Value = '3.430.37'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("E6").Value = '  -4.85%  '
$ws.Range("D7").Value = '3.430.76'
$ws.Range("E7").Value = '  -0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.479'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").Value = '4.017.72'
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.05%  '
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000172'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.434.22'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").Value = '61.266.54'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '398.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.97%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("E26").Value = '  -2.35%  '
$ws.Range("D27").Value = '3.585.93'
$ws.Range("E27").Value = '  -0.32%  '
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.43%  '
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("E32").Value = '  -0.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.38%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.89'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("D36").Value = '3.458.87'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.08%  '
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0787'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.801'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").Value = '2.619.57'
$ws.Range("E47").Value = '  -1.62%  '
$ws.Range("E48").Value = '  -4.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.38'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.07%  '
